$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 964.1818
$ws.Range("I39").Value = 178.44444
$ws.Range("K39").Value = 535.33332
$ws.Range("M39").Value = -239.33332
$ws.Range("H40").Value = 3501.2222
$ws.Range("I40").Value = 3108.2
$ws.Range("J40").Value = 3992.5
$ws.Range("K40").Value = 3108.2
$ws.Range("L40").Value = 3992.5
$ws.Range("M40").Value = -2933.2
$ws.Range("N40").Value = -4342.5
$ws.Range("H51").Value = 23331.334
$ws.Range("J51").Value = 5999.6
$ws.Range("L51").Value = 5999.6
$ws.Range("N51").Value = -6967.6
$ws.Range("H53").Value = 378.6154
$ws.Range("I53").Value = 105.55556
$ws.Range("K53").Value = 105.55556
$ws.Range("M53").Value = 531.44444
$ws.Range("H64").Value = 4926.6
$ws.Range("I64").Value = 4926.6
$ws.Range("K64").Value = 4926.6
$ws.Range("M64").Value = -4678.6
$ws.Range("H67").Value = 4926.6
$ws.Range("I67").Value = 4926.6
$ws.Range("K67").Value = 4926.6
$ws.Range("M67").Value = -4068.6
$ws.Range("H70").Value = 1654.5
$ws.Range("I70").Value = 1490.2
$ws.Range("K70").Value = 4470.6
$ws.Range("M70").Value = -4200.6
$ws.Range("H73").Value = 1654.5
$ws.Range("I73").Value = 1490.2
$ws.Range("K73").Value = 4470.6
$ws.Range("M73").Value = -3534.6
$ws.Range("H92").Value = 383.42105
$ws.Range("I92").Value = 98.94118
$ws.Range("K92").Value = 98.94118
$ws.Range("M92").Value = 1149.05882
$ws.Range("H94").Value = 15015
$ws.Range("I94").Value = 16818
$ws.Range("J94").Value = 6000
$ws.Range("K94").Value = 16818
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = -16367
$ws.Range("N94").Value = -6902
$ws.Range("H98").Value = 1145.3103
$ws.Range("I98").Value = 1152.037
$ws.Range("J98").Value = 1054.5
$ws.Range("K98").Value = 1152.037
$ws.Range("L98").Value = 1054.5
$ws.Range("M98").Value = 345.963
$ws.Range("N98").Value = -4050.5
$ws.Range("H99").Value = 62504090
$ws.Range("I99").Value = 800.1539
$ws.Range("K99").Value = 2400.4617
$ws.Range("M99").Value = -902.4616999999998
$ws.Range("H100").Value = 1574.1765
$ws.Range("I100").Value = 1047.9
$ws.Range("J100").Value = 2326
$ws.Range("K100").Value = 1047.9
$ws.Range("L100").Value = 2326
$ws.Range("M100").Value = -506.9000000000001
$ws.Range("N100").Value = -3408
$ws.Range("H122").Value = 1145.3103
$ws.Range("I122").Value = 1152.037
$ws.Range("J122").Value = 1054.5
$ws.Range("K122").Value = 3456.111
$ws.Range("L122").Value = 3163.5
$ws.Range("M122").Value = -1006.111
$ws.Range("N122").Value = -8063.5
$ws.Range("H129").Value = 25653786
$ws.Range("I129").Value = 4722.25
$ws.Range("J129").Value = 37053370
$ws.Range("K129").Value = 14166.75
$ws.Range("L129").Value = 111160110
$ws.Range("M129").Value = -9166.75
$ws.Range("N129").Value = -111170110

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2237.1428
$ws.Range("I63").Value = 2599.5715
$ws.Range("K63").Value = 2599.5715
$ws.Range("M63").Value = -1913.5715
$ws.Range("H66").Value = 2237.1428
$ws.Range("I66").Value = 2599.5715
$ws.Range("K66").Value = 12997.8575
$ws.Range("M66").Value = -9565.8575
$ws.Range("H97").Value = 841.67566
$ws.Range("I97").Value = 831.03845
$ws.Range("K97").Value = 831.03845
$ws.Range("M97").Value = -335.03845
$ws.Range("H102").Value = 5771.909
$ws.Range("I102").Value = 5525.421
$ws.Range("K102").Value = 5525.421
$ws.Range("M102").Value = -3903.421
$ws.Range("H110").Value = 4404.0586
$ws.Range("I110").Value = 3002.1052
$ws.Range("J110").Value = 6179.8667
$ws.Range("K110").Value = 3002.1052
$ws.Range("L110").Value = 6179.8667
$ws.Range("M110").Value = -957.1052
$ws.Range("N110").Value = -10269.8667
$ws.Range("H122").Value = 3203.85
$ws.Range("I122").Value = 2981.889
$ws.Range("J122").Value = 3664.8462
$ws.Range("K122").Value = 8945.667000000001
$ws.Range("L122").Value = 10994.5386
$ws.Range("M122").Value = -6495.667000000001
$ws.Range("N122").Value = -15894.5386
$ws.Range("H132").Value = 23329.318
$ws.Range("I132").Value = 1689.2693
$ws.Range("J132").Value = 54587.168
$ws.Range("K132").Value = 5067.8079
$ws.Range("L132").Value = 163761.504
$ws.Range("M132").Value = -2537.8079
$ws.Range("N132").Value = -168821.504

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7390.5
$ws.Range("I31").Value = 2118.4
$ws.Range("J31").Value = 11156.286
$ws.Range("K31").Value = 2118.4
$ws.Range("L31").Value = 11156.286
$ws.Range("M31").Value = -1823.4
$ws.Range("N31").Value = -11746.286
$ws.Range("H34").Value = 7390.5
$ws.Range("I34").Value = 2118.4
$ws.Range("J34").Value = 11156.286
$ws.Range("K34").Value = 2118.4
$ws.Range("L34").Value = 11156.286
$ws.Range("M34").Value = -1916.4
$ws.Range("N34").Value = -11560.286
$ws.Range("H132").Value = 5368.5435
$ws.Range("I132").Value = 4658.825
$ws.Range("K132").Value = 13976.475
$ws.Range("M132").Value = -11446.475
$ws.Range("H134").Value = 12191.104
$ws.Range("I134").Value = 12786.429
$ws.Range("K134").Value = 38359.287
$ws.Range("M134").Value = -35824.287

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3115.1667
$ws.Range("I81").Value = 2283
$ws.Range("K81").Value = 6849
$ws.Range("M81").Value = -5726
$ws.Range("H84").Value = 3115.1667
$ws.Range("I84").Value = 2283
$ws.Range("K84").Value = 20547
$ws.Range("M84").Value = -14931

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6429.8276
$ws.Range("I122").Value = 5950.05
$ws.Range("K122").Value = 17850.15
$ws.Range("M122").Value = -15400.15
$ws.Range("H132").Value = 1906.9259
$ws.Range("I132").Value = 1820.3889
$ws.Range("J132").Value = 2080
$ws.Range("K132").Value = 5461.1667
$ws.Range("L132").Value = 6240
$ws.Range("M132").Value = -2931.1667
$ws.Range("N132").Value = -11300

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3942.7878
$ws.Range("I122").Value = 2794.7222
$ws.Range("K122").Value = 8384.1666
$ws.Range("M122").Value = -5934.1666
